$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Row 24: task's status reverts from "Hecho" to "En proceso"; its estimated-hours
# and day-10 (AI) consumption entries are cleared (task no longer finished).
$ws.Range("F24").Value = "En proceso"
$ws.Range("G24").ClearContents()
$ws.Range("AI24").ClearContents()

# Row 26: status reverts from "Hecho" to "En proceso".
$ws.Range("F26").Value = "En proceso"

# Row 27: status reverts from "Hecho" to "Por iniciar"; day-11 (AL) consumption cleared.
$ws.Range("F27").Value = "Por iniciar"
$ws.Range("AL27").ClearContents()

# Row 28: status reverts from "Hecho" to "Por iniciar".
$ws.Range("F28").Value = "Por iniciar"

# Row 29: status reverts from "Hecho" to "Por iniciar"; day-11 (AL) consumption cleared.
$ws.Range("F29").Value = "Por iniciar"
$ws.Range("AL29").ClearContents()

# Row 30: status reverts from "En proceso" to "Por iniciar"; consumption moves back
# from day-11 (AL) to day-10 (AI).
$ws.Range("F30").Value = "Por iniciar"
$ws.Range("AI30").Value = 1

# Restore the active-cell selection on the "Casos de Uso" sheet to AI31.
$ws.Range("AI31").Select() | Out-Null
